{"js": "const replacements = [\n  [\"2025-06-30 Monday\", \"2025-07-01 Tuesday\"],\n  [\"44\u00d795=\", \"46\u00d720=\"],\n  [\"75\u00d733=\", \"19\u00d748=\"],\n  [\"78\u00d781=\", \"34\u00d728=\"],\n  [\"39\u00d729=\", \"37\u00d768=\"],\n  [\"27\u00d714=\", \"63\u00d773=\"],\n  [\"18\u00d728=\", \"37\u00d764=\"],\n  [\"62\u00d795=\", \"97\u00d717=\"],\n  [\"46\u00d735=\", \"75\u00d735=\"],\n  [\"65\u00d715=\", \"89\u00d721=\"],\n  [\"15\u00d736=\", \"83\u00d769=\"],\n  [\"32\u00d717=\", \"40\u00d777=\"],\n  [\"52\u00d778=\", \"84\u00d714=\"],\n  [\"37\u00d766=\", \"67\u00d795=\"],\n  [\"11\u00d776=\", \"31\u00d774=\"],\n  [\"70\u00d750=\", \"72\u00d756=\"],\n  [\"82\u00d725=\", \"11\u00d757=\"],\n  [\"74\u00d734=\", \"24\u00d765=\"],\n  [\"15\u00d720=\", \"69\u00d751=\"],\n  [\"58\u00d798=\", \"50\u00d714=\"],\n  [\"81\u00d752=\", \"64\u00d771=\"],\n  [\"84\u00d771=\", \"68\u00d794=\"],\n  [\"77\u00d762=\", \"87\u00d789=\"],\n  [\"76\u00d761=\", \"25\u00d790=\"],\n  [\"79\u00d719=\", \"63\u00d779=\"],\n  [\"73\u00d784=\", \"70\u00d780=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"2025-06-30 Monday\", \"2025-07-01 Tuesday\")\n  ,@(\"44\u00d795=\", \"46\u00d720=\")\n  ,@(\"75\u00d733=\", \"19\u00d748=\")\n  ,@(\"78\u00d781=\", \"34\u00d728=\")\n  ,@(\"39\u00d729=\", \"37\u00d768=\")\n  ,@(\"27\u00d714=\", \"63\u00d773=\")\n  ,@(\"18\u00d728=\", \"37\u00d764=\")\n  ,@(\"62\u00d795=\", \"97\u00d717=\")\n  ,@(\"46\u00d735=\", \"75\u00d735=\")\n  ,@(\"65\u00d715=\", \"89\u00d721=\")\n  ,@(\"15\u00d736=\", \"83\u00d769=\")\n  ,@(\"32\u00d717=\", \"40\u00d777=\")\n  ,@(\"52\u00d778=\", \"84\u00d714=\")\n  ,@(\"37\u00d766=\", \"67\u00d795=\")\n  ,@(\"11\u00d776=\", \"31\u00d774=\")\n  ,@(\"70\u00d750=\", \"72\u00d756=\")\n  ,@(\"82\u00d725=\", \"11\u00d757=\")\n  ,@(\"74\u00d734=\", \"24\u00d765=\")\n  ,@(\"15\u00d720=\", \"69\u00d751=\")\n  ,@(\"58\u00d798=\", \"50\u00d714=\")\n  ,@(\"81\u00d752=\", \"64\u00d771=\")\n  ,@(\"84\u00d771=\", \"68\u00d794=\")\n  ,@(\"77\u00d762=\", \"87\u00d789=\")\n  ,@(\"76\u00d761=\", \"25\u00d790=\")\n  ,@(\"79\u00d719=\", \"63\u00d779=\")\n  ,@(\"73\u00d784=\", \"70\u00d780=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    Write-Output \"WARNING: no match for $oldText\"\n  }\n}"}
